$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the cell values can be updated.
$ws.Unprotect("D382")

# Update the confidential disclaimer date (2021-03-18 -> 2021-03-19) in place.
$null = $ws.Range("A12").Replace("2021-03-18", "2021-03-19")

# Update Weight (D) and Percent Change (E) values for rows 2-9.
$ws.Range("D2").Value = 0.1778042623096033
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 0.1775989177131859
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.2238084980452
$ws.Range("E4").Value = 0.000841042893187538

$ws.Range("D5").Value = 0.08014780784577413
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.08019713081256065
# E6 is unchanged (stays 0).

$ws.Range("D7").Value = 0.120559423102475
$ws.Range("E7").Value = -0.0009852216748769127

$ws.Range("D8").Value = 0.139883960171201
$ws.Range("E8").Value = 0.0008347245409014992

$ws.Range("E9").Value = 0.0001862193643982213

# Re-protect the sheet the way it was before.
$ws.Protect("D382")
